$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Reference subst"
$ws.Range("D10").Value = "my-ref-subst"

$ws.Range("D9").Font.Color = $ws.Range("C9").Font.Color

$ws.Range("H8").Select()
